$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hand ")
$ws.Cells.Item(19, 1).Value = "Hand19"
$ws.Cells.Item(19, 2).Value = "Eoin mac Domhnaill Ó Conchubuir"
$ws.Cells.Item(19, 3).Value = "Transcription 14"

$ws.Activate()
$ws.Range("G9").Select()
